$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.04172184405617529
$ws.Range("C2").Value = 0.3048912486333797
$ws.Range("D2").Value = 3.223369029078222
$ws.Range("E2").Value = 0.5333859586016987
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 4.103368080369475
